# Update "Price" (column D) and "Volume(1h)" (column E) figures in the
# cryptos list to the latest scraped values.
#
# For column D values that look like a plain decimal number (e.g. "587.02"),
# a leading apostrophe is used to force Excel to store them as text (matching
# the original cell content, which is a text string such as "134.17" or
# "64.239.15" rather than a numeric value), and the cell style is reset to
# "Normal" afterwards so no stray number formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.252.24"
$ws.Range("D3").Value = "3.488.79"
$ws.Range("E3").Value = "  -1.08%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'587.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("D6").Value = "'134.23"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.54%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  +0.34%  "
$ws.Range("E9").Value = "  -0.15%  "
$ws.Range("D10").Value = "'7.27"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.35%  "
$ws.Range("E11").Value = "  +1.89%  "
$ws.Range("D12").Value = "4.081.67"
$ws.Range("E12").Value = "  -1.09%  "
$ws.Range("E13").Value = "  +1.65%  "
$ws.Range("E14").Value = "  +0.92%  "
$ws.Range("D15").Value = "3.489.53"
$ws.Range("E15").Value = "  -0.77%  "
$ws.Range("D16").Value = "'25.80"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.83%  "
$ws.Range("D17").Value = "64.318.93"
$ws.Range("E17").Value = "  +0.31%  "
$ws.Range("D18").Value = "'9.90"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("E19").Value = "  +2.29%  "
$ws.Range("E20").Value = "  -3.74%  "
$ws.Range("D21").Value = "'394.39"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.57%  "
$ws.Range("D22").Value = "'0.570"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.90%  "
$ws.Range("D23").Value = "3.628.66"
$ws.Range("D24").Value = "'74.58"
$ws.Range("D24").Style = "Normal"
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("E26").Value = "  +0.30%  "
$ws.Range("E27").Value = "  -0.31%  "
$ws.Range("D28").Value = "'0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("D29").Value = "'7.41"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.53%  "
$ws.Range("E30").Value = "  -5.75%  "
$ws.Range("D31").Value = "'2.23"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.27%  "
$ws.Range("E32").Value = "  -0.90%  "
$ws.Range("D33").Value = "3.511.71"
$ws.Range("E33").Value = "  -0.68%  "
$ws.Range("E34").Value = "  +3.67%  "
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("D36").Value = "'23.42"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.61%  "
$ws.Range("D37").Value = "'5.14"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.37%  "
$ws.Range("E38").Value = "  -0.21%  "
$ws.Range("E39").Value = "  -1.43%  "
$ws.Range("D40").Value = "'166.22"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.47%  "
$ws.Range("D41").Value = "'0.0781"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.23%  "
$ws.Range("E42").Value = "  -1.23%  "
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("D44").Value = "'25.31"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.61%  "
$ws.Range("E45").Value = "  -0.70%  "
$ws.Range("E46").Value = "  +2.11%  "
$ws.Range("E47").Value = "  -3.73%  "
$ws.Range("D48").Value = "2.457.39"
$ws.Range("E48").Value = "  -0.35%  "
$ws.Range("E49").Value = "  -1.01%  "
$ws.Range("E50").Value = "  -1.78%  "
$ws.Range("D51").Value = "'0.0260"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.47%  "
